# Apply the "Add files via upload" edit to the STAMP Session ID deck.
#
# Summary of changes:
#   1. Slide 1, shape "Rectangle 1": update the "Updated: March 16, 2020"
#      line to "Updated: March 22, 2020".
#   2. Slide 2: delete shape "Rectangle 6" (the floating Timestamp/Control
#      Code diagram rectangle that duplicated the one already on the
#      slide).
#   3. Slide 2, shape "Rectangle 2": reposition/resize the text box and
#      update its "For a Query:" / "For a Response -> ..." headings.

$p = $ppt.ActivePresentation

# --- Slide 1: bump the "Updated" date ------------------------------------
$slide1 = $p.Slides.Item(1)
$rect1 = $slide1.Shapes.Item(4)   # "Rectangle 1"
$origHeight1 = $rect1.Height      # the autofit box re-lays-out on text edit;
                                   # preserve its original (unchanged) extent
$tr1 = $rect1.TextFrame.TextRange
$dateRun = $tr1.Paragraphs(2).Runs(1)
$dateRun.Text = "Updated: March 22, 2020"
$rect1.Height = $origHeight1

# --- Slide 2: remove the duplicate "Rectangle 6" diagram -----------------
$slide2 = $p.Slides.Item(2)
$rect6 = $slide2.Shapes.Item(4)   # "Rectangle 6"
$rect6.Delete()

# --- Slide 2: resize/move + retitle the "Rectangle 2" text box -----------
$rect2 = $slide2.Shapes.Item(4)   # "Rectangle 2" (now at index 4 post-delete)

$rect2.Left = 102
$rect2.Top = 70.5
$rect2.Width = 408
$rect2.Height = 278.6952755905512

$rect2tr = $rect2.TextFrame.TextRange
$rect2tr.Paragraphs(1).Runs(1).Text = "For a Query: TX Control Code"
$rect2tr.Paragraphs(9).Runs(1).Text = "For a Response -> RX Control Code - Used to return Errors"
